# Progress-Check/Timeline.xlsx — mark several tasks as Complete / Out of
# Scope (clearing their start/end dates) and push out the Phase Vocoder
# task's schedule by a week.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# Tasks that are now fully Complete -> clear start/end dates, tag status.
$completeRows = 9,10,11,12,13
foreach ($r in $completeRows) {
    $ws.Cells.Item($r, 3).Value = "Complete"
    $ws.Cells.Item($r, 5).ClearContents()
    $ws.Cells.Item($r, 6).ClearContents()
}

# Gates/Triggers Component is no longer planned.
$ws.Cells.Item(15, 3).Value = "No In Scope"
$ws.Cells.Item(15, 5).ClearContents()
$ws.Cells.Item(15, 6).ClearContents()

# MIDI Keyboard Support is out of scope entirely.
$ws.Cells.Item(17, 3).Value = "Not in Scope"
$ws.Cells.Item(17, 5).ClearContents()
$ws.Cells.Item(17, 6).ClearContents()

# Basic Phase Vocoder implementation slips a week and grows from 8 to 15
# days.
$ws.Range("E16").Value = 44312
$ws.Range("F16").Value = 44326

# Leave the selection on the project start date, like the author did.
$ws.Activate() | Out-Null
$ws.Range("E3:F3").Select() | Out-Null
